# "semana 34 de 2025"
#
# The sheet is a weekly IRA-UCI surveillance table: column A/B/C identify the
# reporting facility, and columns D..AJ hold one column per epidemiological
# week (week 1 in column D ... week 33 in column AJ). This commit appends the
# data for week 34 in a new column AK (and also backfills a handful of
# previously-missing week values that arrived late for a few facilities).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: label the new column "34" -----------------------------
# Use a leading apostrophe so Excel keeps this numeric-looking label as text,
# matching the existing week-number headers in D1:AJ1.
$ws.Range("AK1").Value = "'34"

# --- Week 34 counts per facility (column AK) ----------------------------
# Facilities that had no case data for week 34 are left blank, exactly like
# the existing rows that already have gaps in earlier week columns.

$ws.Range("AK2").Value  = 0   # INSTITUTO DE DIAGNOSTICO MEDICO SA
$ws.Range("AK3").Value  = 0   # INSTITUTO DE DIAGNOSTICO MEDICO SA
# row 4  (CENTRO MEDICO N° 1) - no week 34 report
$ws.Range("AK5").Value  = 0   # CENTRO MEDICO Y ODONTOLOGICO DE LA CIRCUNVALAR
$ws.Range("AK6").Value  = 5   # CLINICA COMFAMILIAR
$ws.Range("AK7").Value  = 0   # UNIDAD INTERMEDIA DE KENNEDY
$ws.Range("AK8").Value  = 0   # UNIDAD INTERMEDIA DE CUBA
# row 9  (CENTRO DE SALUD BOSTON) - no week 34 report
$ws.Range("AK10").Value = 0   # CENTRO DE SALUD SANTA TERESITA
$ws.Range("AK11").Value = 0   # CENTRO DE SALUD SAN NICOLAS
$ws.Range("AK12").Value = 0   # CENTRO DE SALUD SAN CAMILO
$ws.Range("AK13").Value = 0   # CENTRO DE SALUD VILLASANTANA
$ws.Range("AK14").Value = 0   # CENTRO DE SALUD PERLA DEL OTUN
# row 15 (CENTRO DE SALUD VILLA CONSOTA) - no week 34 report
# row 16 (PUESTO DE SALUD PUERTO CALDAS) - no week 34 report
$ws.Range("AK17").Value = 0   # PUESTO DE SALUD CAIMALITO
$ws.Range("AK18").Value = 0   # PUESTO DE SALUD FONDA CENTRAL
$ws.Range("AK19").Value = 0   # PUESTO DE SALUD CRUCERO DE COMBIA
# row 20 (PUESTO DE SALUD LA BELLA) - no week 34 report
# row 21 (PUESTO DE SALUD LA FLORIDA) - no week 34 report
$ws.Range("AK22").Value = 0   # PUESTO DE SALUD ARABIA
$ws.Range("AK23").Value = 0   # PUESTO DE SALUD ALTAGRACIA
$ws.Range("AK24").Value = 0   # PUESTO DE SALUD MORELIA
$ws.Range("AK25").Value = 0   # HOSPITAL DEL CENTRO
# row 26 (CENTRO DE SALUD CASA DEL ABUELO) - no week 34 report
$ws.Range("AK27").Value = 0   # CENTRO DE SALUD EL REMANSO
$ws.Range("AK28").Value = 3   # CLINICA LOS ROSALES
$ws.Range("AK29").Value = 1   # ONCOLOGOS DEL OCCIDENTE SA
$ws.Range("AK30").Value = 6   # ESE HOSPITAL UNIVERSITARIO SAN JORGE DE PEREIRA
$ws.Range("AK31").Value = 0   # FUNDACIÓN LA LIGA AMA SALVAR VIDAS
# row 32 (COSMITET LTDA ...) - no week 34 report
# row 33 (SOCIEDAD INTEGRAL DE ESPECIALISTAS EN SALUD SAS) - no week 34 report
$ws.Range("AK34").Value = 0   # SOCIEDAD COMERCIALIZADORA DE INSUMOS Y SERVICIOS M
$ws.Range("AK35").Value = 3   # SOCIEDAD COMERCIALIZADORA DE INSUMOS Y SERVICIOS M
$ws.Range("AK36").Value = 0   # CORPORACIÓN MEDICA SALUD PARA LOS COLOMBIANOS - CM
$ws.Range("AK37").Value = 0   # SINERGIA GLOBAL EN SALUD SAS
$ws.Range("AK38").Value = 0   # RED MEDICA VITAL SAS
# row 39 (SALUD PYP SAS) - no week 34 report
$ws.Range("AK40").Value = 0   # IPS CENTRO DE MEDICINA INTEGRATIVA SAS
$ws.Range("AK41").Value = 0   # SERVICIOS DE SALUD IPS SURAMERICANA SA
$ws.Range("AK42").Value = 0   # CAJA COLOMBIANA DE SUBSIDIO FAMILIAR COLSUBSIDIO
$ws.Range("AK43").Value = 0   # CAJA COLOMBIANA DE SUBSIDIO FAMILIAR COLSUBSIDIO
$ws.Range("AK44").Value = 0   # GRUPO EMI
$ws.Range("AK45").Value = 0   # VIRREY SOLIS SA PINARES
$ws.Range("AK46").Value = 0   # VIRREY SOLIS IPS LAGO
$ws.Range("AK47").Value = 0   # VIRREY SOLIS IPS SA ALPES
$ws.Range("AK48").Value = 0   # VIRREY SOLIS IPS SA LA REBECA
$ws.Range("AK49").Value = 0   # CENTRO MEDICO PEREIRA COLSANITAS
$ws.Range("AK50").Value = 0   # NEUROMEDICA SAS
$ws.Range("AK51").Value = 0   # SPORT MEDICAL IPS GUSTAVO PORTELA SAS
$ws.Range("AK52").Value = 0   # COOMEVA EMERGENCIA MEDICA SERVICIO DE AMBULANCIA P
$ws.Range("AK53").Value = 0   # CLINICA LOS NEVADOS SAS
$ws.Range("AK54").Value = 0   # CLINICA CENTRAL DEL EJE SAS
$ws.Range("AK55").Value = 0   # CLINICA MEDICA TURIN SAS
$ws.Range("AK56").Value = 0   # EPMSC PEREIRA
$ws.Range("AK57").Value = 0   # SANIDAD POLICIA NACIONAL RISARALDA
$ws.Range("AK58").Value = 0   # BATALLON SAN MATEO

# --- Late-arriving corrections for earlier weeks ------------------------
# A few facilities sent amended counts together with this week's update.
$ws.Range("AJ28").Value = 9   # CLINICA LOS ROSALES - week 33 (was missing)
$ws.Range("AJ30").Value = 5   # ESE HOSPITAL UNIVERSITARIO SAN JORGE - week 33 (was missing)
$ws.Range("H35").Value  = 1   # SOCIEDAD COMERCIALIZADORA ... - week 5 correction
$ws.Range("K35").Value  = 2   # SOCIEDAD COMERCIALIZADORA ... - week 8 correction
$ws.Range("U35").Value  = 3   # SOCIEDAD COMERCIALIZADORA ... - week 18 correction
